$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("TextBoxes")

# "height" and "width" were redundant text-box properties; drop both
# columns (they shift the old "textBoxType" column left into their place).
[void]$ws.Range("I1:J1").EntireColumn.Delete()

# leave the cursor where the author left it after the edit
[void]$ws.Range("F16").Select()
